# Gaussian Quadrature Scheme export: rename sheet, correct tiny rounding
# differences on row 13, and append a new "HexGrid-60degTilt5degRes" data
# row (row 16) to the Averaged Intensities sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab (workbook.xml <sheet name="...">)
$ws.Name = "CopperA"

# 2. Tiny floating point corrections on row 13
$ws.Range("L13").Value = 0.9927591642822557
$ws.Range("P13").Value = 0.9941924247502772

# 3. Append new row 16 of averaged-intensity data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.067466239315418
$ws.Range("D16").Value = 0.8773520961605883
$ws.Range("E16").Value = 1.070699185305034
$ws.Range("F16").Value = 0.9439715217766051
$ws.Range("G16").Value = 1.067466239315418
$ws.Range("H16").Value = 0.8773520961605883
$ws.Range("I16").Value = 1.055506493621435
$ws.Range("J16").Value = 0.9784606975843092
$ws.Range("K16").Value = 1.003181088641286
$ws.Range("L16").Value = 0.9002881726143893
$ws.Range("M16").Value = 1.067466239315418
$ws.Range("N16").Value = 0.9740256407328109
$ws.Range("O16").Value = 0.9898722606394112
$ws.Range("P16").Value = 0.9871156868773832

# Carry the row-header style (bold, centered, thin border) from A15 to A16,
# matching the formatting of every other row-number cell in column A.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A16").Value = 14
